$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns.Item(28).Insert()
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:AE15"))
